$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 keeps its existing TC_01 / Ford Focus cells (they keep their existing
# style), we only need to fill in the rest of the row.
$ws.Range("C2").Value = "auto"
$ws.Range("D2").Value = 84102
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 2000
$ws.Range("G2").Value = 10000

# Rows 3 and 4 are fully cleared (content + formatting) and retyped with the
# new data so the resulting cells carry no explicit style, matching a fresh
# entry.
$ws.Range("A3:G4").Clear()

$ws.Range("A3").Value = "TC_02"
$ws.Range("B3").Value = "Skoda Octavia"
$ws.Range("C3").Value = "auto"
$ws.Range("D3").Value = 95801
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 2000
$ws.Range("G3").Value = 5000

$ws.Range("A4").Value = "TC_03"
$ws.Range("B4").Value = "Dacia Duster"
$ws.Range("C4").Value = "auto"
$ws.Range("D4").Value = 84102
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 7000
$ws.Range("G4").Value = 13000

# Remove the now-obsolete trailing test case rows.
$ws.Rows("5:7").Delete()

# Restore the print setup (paper size / orientation) for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Match the final selection left on the sheet.
$null = $ws.Range("H5").Select()
